$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header labels in row 2 ---
$ws.Range("A2").Value = "Target"
$ws.Range("B2").Value = "Feature"

# --- Header row 3 gets bigger font + gray fill, row height grows ---
$ws.Range("A3:B3").Font.Size = 16
$ws.Range("A3:B3").Interior.Color = 12566463

# --- Data rows 4-8 in columns A:B get the bigger font (no fill) ---
$ws.Range("A4:B8").Font.Size = 16

# --- Row heights for rows 3-8 grow to fit the larger font ---
$ws.Rows("3:8").RowHeight = 21

# --- Column widths grow to fit new font/content ---
$ws.Columns("A:A").ColumnWidth = 16
$ws.Columns("B:B").ColumnWidth = 14.3333333333333

# --- Move the old F20 formula down into a new row 21 (E21) ---
$ws.Range("E20").Copy()
$ws.Range("E21").PasteSpecial(-4122)
$ws.Range("F20").Clear()
$ws.Range("E21").Formula = "=J14*E17+J15*F17+J16*G17+H17*J17+I17*J18"

# --- New matrix-inverse block (rows 23-25) ---
$ws.Range("E23").Value = 1
$ws.Range("F23").Value = 601
$ws.Range("G23").Value = -53
$ws.Range("F23:G24").Interior.Color = 12566463
$ws.Range("I23").Formula = "=`$F`$23*E25"
$ws.Range("J23").Formula = "=`$G`$23*E25"
$ws.Range("L23").Value = "A"
$ws.Range("M23").Formula = "=F23*E20 + G23*E21"

$ws.Range("E24").Formula = "=E8*F9 - F8*E9"
$ws.Range("F24").Value = -53
$ws.Range("G24").Value = 5
$ws.Range("I24").Formula = "=`$F`$24*E25"
$ws.Range("J24").Formula = "=`$G`$24*E25"
$ws.Range("M24").Formula = "=F24*E20 + G24*E21"

$ws.Range("E25").Formula = "=E23/E24"

# --- Selection / view ---
[void]$ws.Range("O18").Select()

Write-Host "Edit applied"
